$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset demscalar_ippu_* rows to 1 across columns J:AS (cols 10-45)
$ws.Range("J4:AS10").Value = 1
$ws.Range("J12:AS15").Value = 1
$ws.Range("J21:AS22").Value = 1

# Update prodinit_ippu_* rows with newly computed production values
$row96 = @(14282374.9077076, 11886892.519324, 11428581.2771974, 10821851.2150115, 11487972.2999916, 10215273.3000134, 10215270, 10184278.4388639, 10153380.9013672, 10122577.1022564, 10091866.7571437, 10061249.5825041, 10030725.2956726, 10000293.6148419, 9969954.2590595, 9939706.948225491, 9909551.40308959, 9879487.345248761, 9849514.49714458, 9819632.582060721, 9789841.324120341, 9760140.448283611, 9730529.68034507, 9701008.74693121, 9671577.375497861, 9642235.29432773, 9612982.23252786, 9583817.920027141, 9554742.08757383, 9525754.46673304, 9496854.78988428, 9468042.790218979, 9439318.20173802, 9410680.75924929, 9382130.19836523, 9353666.255500371)
for ($i = 0; $i -lt $row96.Length; $i++) {
    $ws.Cells.Item(96, 10 + $i).Value = $row96[$i]
}

$row97 = @(178948.375716022, 177754.925772729, 180314.442487874, 230718.134473586, 236094.868192736, 204269.118490759, 204269.1, 213052.0937904, 222212.731482508, 231767.250696422, 241732.587223997, 252126.405048258, 262967.127653561, 274273.970682, 286066.975993956, 298367.047193142, 311195.986679131, 324576.534293041, 338532.407624863, 353088.34405391, 368270.144596862, 384104.719641174, 400620.136644881, 417845.669887361, 435811.852359239, 454550.529883418, 474094.917563146, 494479.658657211, 515740.88598658, 537916.28598137, 561045.165481645, 585168.521410464, 610329.11344267, 636571.539798236, 663942.316294502, 692489.9587974499)
for ($i = 0; $i -lt $row97.Length; $i++) {
    $ws.Cells.Item(97, 10 + $i).Value = $row97[$i]
}

$row98 = @(213269.959730345, 243307.677450407, 248364.690415284, 304033.027738509, 337324.805659277, 335571.029555977, 335571, 358363.800442098, 382704.743459071, 408698.982668976, 436458.813979972, 466104.16071208, 497763.091669589, 531572.374402205, 567678.066045051, 606236.144289956, 647413.18121385, 691387.06287523, 738347.757787375, 788498.137588128, 842054.8534515731, 899249.2720277359, 960328.474953585, 1025556.32625325, 1095214.6122387, 1169604.25883518, 1249046.63159045, 1333884.92398376, 1424485.64003218, 1521240.17759915, 1624566.51924524, 1734911.03792555, 1852750.42533457, 1978593.75122883, 2112984.66262385, 2256503.73236587)
for ($i = 0; $i -lt $row98.Length; $i++) {
    $ws.Cells.Item(98, 10 + $i).Value = $row98[$i]
}

$row99 = @(32105.947345554, 31208.0322434067, 33143.0079769574, 35817.5456594234, 38647.3460647448, 40022.0889332774, 40022.09, 42821.9433775996, 45817.6680586733, 49022.9667491669, 52452.5007648219, 56121.9570933129, 60048.1201478963, 64248.9485407812, 68743.6572273873, 73552.8053972268, 78698.39051343031, 84203.9489310612, 90094.6635544569, 96397.4790260289, 103141.224973407, 110356.747878672, 118077.052172851, 126337.451201066, 135175.728748857, 144632.311868503, 154750.455795896, 165576.441803754, 177159.788896201, 189553.480313013, 202814.205879573, 217002.621311078, 232183.625657073, 248426.658155358, 265806.015853121, 284401.193448111)
for ($i = 0; $i -lt $row99.Length; $i++) {
    $ws.Cells.Item(99, 10 + $i).Value = $row99[$i]
}

$row100 = @(1680065.35757331, 888570.984653814, 1060271.201573, 1050411.22764009, 956642.3367598359, 858774.114672421, 858774.1, 842148.292679413, 825844.359841372, 809856.070018103, 794177.312382622, 778802.094413129, 763724.539602635, 748938.885212921, 734439.480071982, 720220.7824141131, 706277.35776181, 692603.876848679, 679195.113582556, 666045.9430480601, 653151.33954782, 640506.374681617, 628106.21546272, 615946.122470691, 604021.448039941, 592327.634483368, 580860.212350376, 569614.798718628, 558587.095518856, 547772.887892121, 537168.042578858, 526768.506339127, 516570.304403428, 506569.538953529, 496762.387632687, 487145.102084723)
for ($i = 0; $i -lt $row100.Length; $i++) {
    $ws.Cells.Item(100, 10 + $i).Value = $row100[$i]
}

$row101 = @(1131515.28131103, 1071341.85070564, 1283367.95584282, 1547203.12108103, 1660328.60479078, 1701991.05060139, 1701991, 1813965.98559712, 1933307.87113642, 2060501.32928356, 2196062.919603, 2340543.18641577, 2494528.89467648, 2658645.41295008, 2833559.25316601, 3019980.77746423, 3218667.08312633, 3430425.07730786, 3656114.75405912, 3896652.68694314, 4153015.75143494, 4426245.09221967, 4717450.35150168, 5027814.17549637, 5358597.01740694, 5711142.25639172, 6086881.65331132, 6487341.1654127, 6914147.14356459, 7369032.93721276, 7853845.93387909, 8370555.06179348, 8921258.786128251, 9508193.63130953, 10133743.2640165, 10800448.1737564)
for ($i = 0; $i -lt $row101.Length; $i++) {
    $ws.Cells.Item(101, 10 + $i).Value = $row101[$i]
}

$row102 = @(5693573.41029101, 2978110.69989548, 2086345.8801511, 3230715.29922382, 3785583.16375786, 774225.61666885, 774225.6, 811491.448961632, 850551.017349271, 891490.642371616, 934400.816911608, 979376.389551745, 1026516.77422725, 1075926.16997047, 1127713.79123227, 1181994.10928945, 1238887.10527191, 1298518.53536874, 1361020.2087995, 1426530.27916511, 1495193.54982234, 1567161.79395698, 1642594.09006298, 1721657.17356933, 1804525.80539182, 1891383.15822425, 1982421.22142309, 2077841.22538044, 2177854.08632341, 2282680.87252293, 2392553.29294283, 2507714.20940886, 2628418.17343016, 2754931.98885949, 2887535.30163603, 3026521.21791438)
for ($i = 0; $i -lt $row102.Length; $i++) {
    $ws.Cells.Item(102, 10 + $i).Value = $row102[$i]
}

$row103 = @(638416.0051241081, 646797.418758036, 716952.084571566, 665744.559176599, 712445.121981666, 710539.46105489, 710539.5, 746749.054770675, 784803.871988674, 824797.987427734, 866830.228985743, 911004.460894505, 957429.840374589, 1006221.08736947, 1057498.76802553, 1111389.59261831, 1168026.72866137, 1227550.12997133, 1290106.88250225, 1355851.56780404, 1424946.64500285, 1497562.85224752, 1573879.62861391, 1654085.55750983, 1738378.83267614, 1826967.74793556, 1920071.21189939, 2017919.28890407, 2120753.76751416, 2228828.75799663, 2342411.32024283, 2461782.12368977, 2587236.14087138, 2719083.37631368, 2857649.63257481, 3003277.31532309)
for ($i = 0; $i -lt $row103.Length; $i++) {
    $ws.Cells.Item(103, 10 + $i).Value = $row103[$i]
}

$row104 = @(444829.300044985, 464726.510804118, 495133.552705071, 534335.928743316, 545622.112267746, 558346.3758981921, 558346.4, 592916.448176721, 629626.902794571, 668610.287236502, 710007.330078038, 753967.473112224, 800649.410828972, 850221.663296296, 902863.1845115219, 958764.0084186, 1018125.93492361, 1081163.25838496, 1148103.54120807, 1219188.43533729, 1294674.5546104, 1374834.40112514, 1459957.34896166, 1550350.68881227, 1646340.73728943, 1748274.01491672, 1856518.49705501, 1971464.94228, 2093528.30300613, 2223149.22344951, 2360795.63033723, 2506964.42210548, 2662183.26268457, 2827012.4863462, 3002047.12048942, 3187919.03366752)
for ($i = 0; $i -lt $row104.Length; $i++) {
    $ws.Cells.Item(104, 10 + $i).Value = $row104[$i]
}

$row111 = @(10718.1652403502, 10116.9427136264, 9452.581899016281, 9076.28496846161, 9512.48698132631, 8252.990289068461, 8252.99, 8425.81475686404, 8602.258613785771, 8782.397357854181, 8966.30836320552, 9154.07062425746, 9345.764789639161, 9541.473196831939, 9741.27990753524, 9945.270743773301, 10153.5333247579, 10366.1571045228, 10583.2334103469, 10804.8554819814, 11031.1185116985, 11262.119685179, 11497.958223256, 11738.735424533, 11984.5547088938, 12235.5216619239, 12491.7440802628, 12753.3320179047, 13020.3978334699, 13293.0562384652, 13571.4243465561, 13855.6217238693, 14145.77044035, 14441.9951221935, 14744.4230053758, 15053.1839903043)
for ($i = 0; $i -lt $row111.Length; $i++) {
    $ws.Cells.Item(111, 10 + $i).Value = $row111[$i]
}

$row112 = @(155301.990365553, 151803.77697392, 153900.077552143, 144916.430808095, 151759.735813147, 138331.46203274, 138331.5, 140235.858757676, 142166.434120232, 144123.587002001, 146107.683285891, 148119.093891785, 150158.194845883, 152225.367351, 154320.997857828, 156445.478137185, 158599.205353252, 160782.582137823, 162996.016665576, 165239.92273038, 167514.71982265, 169820.833207771, 172158.694005601, 174528.739271066, 176931.412075864, 179367.161591301, 181836.443172256, 184339.718442312, 186877.455380053, 189450.128406552, 192058.218474063, 194702.21315593, 197382.606737742, 200099.900309733, 202854.601860464, 205647.226371786)
for ($i = 0; $i -lt $row112.Length; $i++) {
    $ws.Cells.Item(112, 10 + $i).Value = $row112[$i]
}

